# Atualização de bases das ligas, do dia: 26-02-2024 às 22:04
#
# This edit permutes the "match data" (everything except the id in column A,
# the Div/Div-Original-Name columns C/D and the Date column E) among a small
# set of rows (130-137, 139-140, 143, 145) on the single worksheet.
# Column B (id/match-code) and columns F..AC (teams, score, odds, etc.) move
# together as a unit from one row to another; A, C, D and E stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that travel together when rows are permuted: B, then F through AC.
$cols = @(2,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)

function Get-RowData($row) {
    $data = @{}
    foreach ($c in $cols) {
        $data[$c] = $ws.Cells.Item($row, $c).Value2
    }
    return $data
}

function Set-RowData($row, $data) {
    foreach ($c in $cols) {
        $ws.Cells.Item($row, $c).Value = $data[$c]
    }
}

# Snapshot the "before" data for every row that participates in the permutation.
$rows = @(130,131,132,133,134,135,136,137,139,140,143,145)
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = Get-RowData $r
}

# new_row -> source_row (source row's OLD data lands on new_row)
$mapping = @{
    130 = 133
    131 = 132
    132 = 131
    133 = 130
    134 = 135
    135 = 136
    136 = 137
    137 = 134
    139 = 140
    140 = 139
    143 = 145
    145 = 143
}

foreach ($r in $rows) {
    $src = $mapping[$r]
    Set-RowData $r $snapshot[$src]
}
